$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 14663
$ws.Range("I21").Value = 4126
$ws.Range("J21").Value = 25200
$ws.Range("K21").Value = 4126
$ws.Range("L21").Value = 25200
$ws.Range("M21").Value = -3658
$ws.Range("N21").Value = -26136
$ws.Range("H23").Value = 14663
$ws.Range("I23").Value = 4126
$ws.Range("J23").Value = 25200
$ws.Range("K23").Value = 4126
$ws.Range("L23").Value = 25200
$ws.Range("M23").Value = -3892
$ws.Range("N23").Value = -25668
$ws.Range("H92").Value = 372.73334
$ws.Range("I92").Value = 353.23077
$ws.Range("K92").Value = 353.23077
$ws.Range("M92").Value = 894.76923
$ws.Range("H97").Value = 999.3333
$ws.Range("J97").Value = 999.3333
$ws.Range("L97").Value = 2997.9999
$ws.Range("N97").Value = -3989.9999
$ws.Range("H107").Value = 166.33333
$ws.Range("I107").Value = 166.33333
$ws.Range("K107").Value = 166.33333
$ws.Range("M107").Value = 1753.66667
$ws.Range("H112").Value = 1315.174
$ws.Range("J112").Value = 1347.3684
$ws.Range("L112").Value = 4042.1052
$ws.Range("N112").Value = -6258.1052
$ws.Range("H133").Value = 73788.46000000001
$ws.Range("J133").Value = 73788.46000000001
$ws.Range("L133").Value = 73788.46000000001
$ws.Range("N133").Value = -83908.46000000001
$ws.Range("H136").Value = 75776.664
$ws.Range("J136").Value = 75776.664
$ws.Range("L136").Value = 75776.664
$ws.Range("N136").Value = -85976.664
$ws.Range("H137").Value = 540164.9
$ws.Range("I137").Value = 2415.1765
$ws.Range("K137").Value = 7245.529500000001
$ws.Range("M137").Value = -4695.529500000001
$ws.Range("H140").Value = 76467.57000000001
$ws.Range("J140").Value = 80760.664
$ws.Range("L140").Value = 80760.664
$ws.Range("N140").Value = -91120.664
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6115.2266
$ws.Range("I32").Value = 2383.1396
$ws.Range("K32").Value = 2383.1396
$ws.Range("M32").Value = -2096.1396
$ws.Range("H74").Value = 79348.84
$ws.Range("I74").Value = 112242.89
$ws.Range("J74").Value = 5337.25
$ws.Range("K74").Value = 112242.89
$ws.Range("L74").Value = 5337.25
$ws.Range("M74").Value = -111368.89
$ws.Range("N74").Value = -7085.25
$ws.Range("H77").Value = 79348.84
$ws.Range("I77").Value = 112242.89
$ws.Range("J77").Value = 5337.25
$ws.Range("K77").Value = 561214.45
$ws.Range("L77").Value = 26686.25
$ws.Range("M77").Value = -556846.45
$ws.Range("N77").Value = -35422.25
$ws.Range("H111").Value = 35000
$ws.Range("J111").Value = 35000
$ws.Range("L111").Value = 35000
$ws.Range("N111").Value = -43180
$ws.Range("H132").Value = 1972.96
$ws.Range("I132").Value = 1908.3043
$ws.Range("J132").Value = 2716.5
$ws.Range("K132").Value = 5724.9129
$ws.Range("L132").Value = 8149.5
$ws.Range("M132").Value = -3194.9129
$ws.Range("N132").Value = -13209.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4486.1875
$ws.Range("I86").Value = 3693.75
$ws.Range("J86").Value = 5278.625
$ws.Range("K86").Value = 3693.75
$ws.Range("L86").Value = 5278.625
$ws.Range("M86").Value = -2570.75
$ws.Range("N86").Value = -7524.625
$ws.Range("H89").Value = 4486.1875
$ws.Range("I89").Value = 3693.75
$ws.Range("J89").Value = 5278.625
$ws.Range("K89").Value = 18468.75
$ws.Range("L89").Value = 26393.125
$ws.Range("M89").Value = -12852.75
$ws.Range("N89").Value = -37625.125
$ws.Range("H105").Value = 59158.668
$ws.Range("I105").Value = 74346.86
$ws.Range("J105").Value = 6000
$ws.Range("K105").Value = 74346.86
$ws.Range("L105").Value = 6000
$ws.Range("M105").Value = -72599.86
$ws.Range("N105").Value = -9494
$ws.Range("H132").Value = 29910.666
$ws.Range("J132").Value = 29910.666
$ws.Range("L132").Value = 29910.666
$ws.Range("N132").Value = -40030.666
$ws.Range("H140").Value = 43499.06
$ws.Range("J140").Value = 43499.06
$ws.Range("L140").Value = 43499.06
$ws.Range("N140").Value = -53859.06
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2578.7036
$ws.Range("I31").Value = 1613
$ws.Range("J31").Value = 5337.857
$ws.Range("K31").Value = 1613
$ws.Range("L31").Value = 5337.857
$ws.Range("M31").Value = -1318
$ws.Range("N31").Value = -5927.857
$ws.Range("H34").Value = 2578.7036
$ws.Range("I34").Value = 1613
$ws.Range("J34").Value = 5337.857
$ws.Range("K34").Value = 1613
$ws.Range("L34").Value = 5337.857
$ws.Range("M34").Value = -1411
$ws.Range("N34").Value = -5741.857
$ws.Range("H117").Value = 36648.8
$ws.Range("J117").Value = 36648.8
$ws.Range("L117").Value = 36648.8
$ws.Range("N117").Value = -45826.8
$ws.Range("H122").Value = 4227.75
$ws.Range("I122").Value = 3237.8
$ws.Range("J122").Value = 4934.857
$ws.Range("K122").Value = 9713.400000000001
$ws.Range("L122").Value = 14804.571
$ws.Range("M122").Value = -7263.400000000001
$ws.Range("N122").Value = -19704.571
$ws.Range("H132").Value = 1896441.1
$ws.Range("I132").Value = 2068549.4
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 6205648.199999999
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -6203118.199999999
$ws.Range("N132").Value = -14810
$ws.Range("H134").Value = 2338943.2
$ws.Range("I134").Value = 2859758.5
$ws.Range("J134").Value = 168878.67
$ws.Range("K134").Value = 8579275.5
$ws.Range("L134").Value = 506636.01
$ws.Range("M134").Value = -8576740.5
$ws.Range("N134").Value = -511706.01
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12091048
$ws.Range("I4").Value = 10000153
$ws.Range("K4").Value = 30000459
$ws.Range("M4").Value = -30000347
$ws.Range("H68").Value = 202499.4
$ws.Range("J68").Value = 252624.25
$ws.Range("L68").Value = 757872.75
$ws.Range("N68").Value = -759494.75
$ws.Range("H71").Value = 202499.4
$ws.Range("J71").Value = 252624.25
$ws.Range("L71").Value = 2273618.25
$ws.Range("N71").Value = -2281730.25
$ws.Range("H97").Value = 112.46667
$ws.Range("I97").Value = 110.125
$ws.Range("K97").Value = 330.375
$ws.Range("M97").Value = 165.625
$ws.Range("H125").Value = 7999
$ws.Range("J125").Value = 10000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840
$ws.Range("H132").Value = 3465.52
$ws.Range("I132").Value = 780.44446
$ws.Range("J132").Value = 4975.875
$ws.Range("K132").Value = 7024.00014
$ws.Range("L132").Value = 44782.875
$ws.Range("M132").Value = -4494.00014
$ws.Range("N132").Value = -49842.875
$ws.Range("H133").Value = 9975
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 9975
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 29925
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -40045
$ws.Range("H138").Value = 7775.023
$ws.Range("I138").Value = 7963.8
$ws.Range("K138").Value = 23891.4
$ws.Range("M138").Value = -18751.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3000.75
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3001
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3001
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -4997
$ws.Range("H83").Value = 3000.75
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3001
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 15005
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -24989
$ws.Range("H140").Value = 91885.664
$ws.Range("J140").Value = 92058.875
$ws.Range("L140").Value = 92058.875
$ws.Range("N140").Value = -102418.875
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1141.8889
$ws.Range("I113").Value = 468.14285
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 1404.42855
$ws.Range("L113").Value = 10500
$ws.Range("M113").Value = 765.5714499999999
$ws.Range("N113").Value = -14840
$ws.Range("H122").Value = 1991.5834
$ws.Range("I122").Value = 1775
$ws.Range("K122").Value = 5325
$ws.Range("M122").Value = -2875
$ws.Range("H132").Value = 2301
$ws.Range("I132").Value = 1837.0714
$ws.Range("J132").Value = 3600
$ws.Range("K132").Value = 5511.2142
$ws.Range("L132").Value = 10800
$ws.Range("M132").Value = -2981.2142
$ws.Range("N132").Value = -15860
